$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "30.262.05"
$ws.Cells.Item(2, 5).Value = "  +1.82%  "
$ws.Cells.Item(3, 4).Value = "2.090.22"
$ws.Cells.Item(3, 5).Value = "  -0.42%  "
$ws.Cells.Item(4, 4).Value = "1.002"
$ws.Cells.Item(4, 5).Value = "  -0.72%  "
$ws.Cells.Item(5, 4).Value = "341.26"
$ws.Cells.Item(5, 5).Value = "  -0.79%  "
$ws.Cells.Item(6, 5).Value = "  -0.63%  "
$ws.Cells.Item(7, 4).Value = "0.5309"
$ws.Cells.Item(7, 5).Value = "  +2.26%  "
$ws.Cells.Item(8, 4).Value = "0.4376"
$ws.Cells.Item(8, 5).Value = "  -0.32%  "
$ws.Cells.Item(9, 4).Value = "54.38"
$ws.Cells.Item(9, 5).Value = "  +1.15%  "
$ws.Cells.Item(10, 4).Value = "0.09369"
$ws.Cells.Item(10, 5).Value = "  +0.98%  "
$ws.Cells.Item(11, 5).Value = "  +0.49%  "
$ws.Cells.Item(12, 5).Value = "  -0.17%  "
$ws.Cells.Item(13, 4).Value = "8.557"
$ws.Cells.Item(13, 5).Value = "  +4.70%  "
$ws.Cells.Item(14, 5).Value = "  +0.95%  "
$ws.Cells.Item(15, 4).Value = "2.019.04"
$ws.Cells.Item(15, 5).Value = "  -5.42%  "
$ws.Cells.Item(16, 5).Value = "  -1.37%  "
$ws.Cells.Item(17, 4).Value = "0.00001157"
$ws.Cells.Item(17, 5).Value = "  +0.20%  "
$ws.Cells.Item(18, 4).Value = "1.002"
$ws.Cells.Item(18, 5).Value = "  -0.67%  "
$ws.Cells.Item(19, 5).Value = "  +0.26%  "
$ws.Cells.Item(20, 4).Value = "0.06712"
$ws.Cells.Item(20, 5).Value = "  +0.67%  "
$ws.Cells.Item(21, 4).Value = "6.339"
$ws.Cells.Item(21, 5).Value = "  +1.97%  "
$ws.Cells.Item(22, 4).Value = "1.000"
$ws.Cells.Item(22, 5).Value = "  -0.76%  "
$ws.Cells.Item(23, 4).Value = "30.244.42"
$ws.Cells.Item(23, 5).Value = "  +1.61%  "
$ws.Cells.Item(24, 4).Value = "12.49"
$ws.Cells.Item(24, 5).Value = "  -0.49%  "
$ws.Cells.Item(25, 5).Value = "  +0.41%  "
$ws.Cells.Item(26, 2).Value = "EthereumClassic"
$ws.Cells.Item(26, 3).Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Cells.Item(26, 4).Value = "21.80"
$ws.Cells.Item(26, 5).Value = "  -0.59%  "
$ws.Cells.Item(27, 2).Value = "InternetComputer(DFINITY)"
$ws.Cells.Item(27, 3).Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Cells.Item(27, 4).Value = "6.890"
$ws.Cells.Item(27, 5).Value = "  +8.05%  "
$ws.Cells.Item(28, 2).Value = "Monero"
$ws.Cells.Item(28, 3).Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Cells.Item(28, 4).Value = "162.58"
$ws.Cells.Item(28, 5).Value = "  +0.19%  "
$ws.Cells.Item(29, 2).Value = "LidoDAOToken"
$ws.Cells.Item(29, 3).Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Cells.Item(29, 4).Value = "2.497"
$ws.Cells.Item(29, 5).Value = "  +0.15%  "
$ws.Cells.Item(30, 2).Value = "BitcoinCash"
$ws.Cells.Item(30, 3).Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Cells.Item(30, 4).Value = "133.67"
$ws.Cells.Item(30, 5).Value = "  +0.01%  "
$ws.Cells.Item(31, 2).Value = "ImmutableX"
$ws.Cells.Item(31, 3).Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Cells.Item(31, 4).Value = "1.128"
$ws.Cells.Item(31, 5).Value = "  -0.33%  "
$ws.Cells.Item(32, 2).Value = "Stellar"
$ws.Cells.Item(32, 3).Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Cells.Item(32, 4).Value = "0.1052"
$ws.Cells.Item(32, 5).Value = "  +0.06%  "
$ws.Cells.Item(33, 2).Value = "ARBITRUM"
$ws.Cells.Item(33, 3).Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Cells.Item(33, 4).Value = "1.665"
$ws.Cells.Item(33, 5).Value = "  -2.24%  "
$ws.Cells.Item(34, 2).Value = "Filecoin"
$ws.Cells.Item(34, 3).Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Cells.Item(34, 4).Value = "6.251"
$ws.Cells.Item(34, 5).Value = "  +0.85%  "
$ws.Cells.Item(35, 2).Value = "HuobiToken"
$ws.Cells.Item(35, 3).Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Cells.Item(35, 4).Value = "3.913"
$ws.Cells.Item(35, 5).Value = "  -1.01%  "
$ws.Cells.Item(36, 2).Value = "FraxShare"
$ws.Cells.Item(36, 3).Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Cells.Item(36, 4).Value = "10.08"
$ws.Cells.Item(36, 5).Value = "  -3.42%  "
$ws.Cells.Item(37, 2).Value = "VeChain"
$ws.Cells.Item(37, 3).Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Cells.Item(37, 4).Value = "0.02614"
$ws.Cells.Item(37, 5).Value = "  +1.23%  "
$ws.Cells.Item(38, 2).Value = "Hedera"
$ws.Cells.Item(38, 3).Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Cells.Item(38, 4).Value = "0.06762"
$ws.Cells.Item(38, 5).Value = "  +0.60%  "
$ws.Cells.Item(39, 2).Value = "Aptos"
$ws.Cells.Item(39, 3).Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Cells.Item(39, 4).Value = "12.55"
$ws.Cells.Item(39, 5).Value = "  +0.23%  "
$ws.Cells.Item(40, 2).Value = "TheSandbox"
$ws.Cells.Item(40, 3).Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Cells.Item(40, 4).Value = "0.6936"
$ws.Cells.Item(40, 5).Value = "  -0.85%  "
$ws.Cells.Item(41, 2).Value = "TrustWalletToken"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Cells.Item(41, 4).Value = "1.340"
$ws.Cells.Item(41, 5).Value = "  +0.29%  "
$ws.Cells.Item(42, 2).Value = "Algorand"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Cells.Item(42, 4).Value = "0.2211"
$ws.Cells.Item(42, 5).Value = "  -0.24%  "
$ws.Cells.Item(43, 2).Value = "Decentraland"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Cells.Item(43, 4).Value = "0.6794"
$ws.Cells.Item(44, 2).Value = "NEARProtocol"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Cells.Item(44, 4).Value = "2.362"
$ws.Cells.Item(44, 5).Value = "  +0.88%  "
$ws.Cells.Item(45, 4).Value = "14.22"
$ws.Cells.Item(45, 5).Value = "  -0.81%  "
$ws.Cells.Item(46, 2).Value = "Frax"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Cells.Item(46, 4).Value = "1.001"
$ws.Cells.Item(46, 5).Value = "  -0.57%  "
$ws.Cells.Item(47, 2).Value = "WEMIXTOKEN"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Cells.Item(47, 4).Value = "1.279"
$ws.Cells.Item(47, 5).Value = "  +6.67%  "
$ws.Cells.Item(48, 2).Value = "PancakeSwap"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Cells.Item(48, 4).Value = "3.632"
$ws.Cells.Item(48, 5).Value = "  +0.17%  "
$ws.Cells.Item(49, 2).Value = "BabyDogeCoin"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Cells.Item(49, 4).Value = "0.00000000348"
$ws.Cells.Item(49, 5).Value = "  -2.65%  "
$ws.Cells.Item(50, 2).Value = "ThetaToken"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Cells.Item(50, 4).Value = "1.205"
$ws.Cells.Item(50, 5).Value = "  +3.49%  "
$ws.Cells.Item(51, 2).Value = "EOS"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Cells.Item(51, 4).Value = "1.212"
$ws.Cells.Item(51, 5).Value = "  -0.51%  "
